$wb = $excel.ActiveWorkbook

function Set-DateText {
    param($ws, $row, $col, $text)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---- INDI: append forecast rows 122-137 (2026Q1 .. 2029Q4) ----
$wsIndi = $wb.Worksheets.Item("INDI")
Set-DateText $wsIndi 122 1 "2026-01-01"
$wsIndi.Cells.Item(122, 2).Value = 1235.04586404976
$wsIndi.Cells.Item(122, 3).Value = 169.26139261989
Set-DateText $wsIndi 123 1 "2026-04-01"
$wsIndi.Cells.Item(123, 2).Value = 1239.5218052771
$wsIndi.Cells.Item(123, 3).Value = 169.874814410511
Set-DateText $wsIndi 124 1 "2026-07-01"
$wsIndi.Cells.Item(124, 2).Value = 1254.5768609402
$wsIndi.Cells.Item(124, 3).Value = 171.938089760587
Set-DateText $wsIndi 125 1 "2026-10-01"
$wsIndi.Cells.Item(125, 2).Value = 1268.66823646484
$wsIndi.Cells.Item(125, 3).Value = 173.869294029721
Set-DateText $wsIndi 126 1 "2027-01-01"
$wsIndi.Cells.Item(126, 2).Value = 1268.65365688043
$wsIndi.Cells.Item(126, 3).Value = 173.867295917073
Set-DateText $wsIndi 127 1 "2027-04-01"
$wsIndi.Cells.Item(127, 2).Value = 1274.59931882781
$wsIndi.Cells.Item(127, 3).Value = 174.682140977127
Set-DateText $wsIndi 128 1 "2027-07-01"
$wsIndi.Cells.Item(128, 2).Value = 1286.8865547456
$wsIndi.Cells.Item(128, 3).Value = 176.366090313288
Set-DateText $wsIndi 129 1 "2027-10-01"
$wsIndi.Cells.Item(129, 2).Value = 1298.88413915914
$wsIndi.Cells.Item(129, 3).Value = 178.01034329613
Set-DateText $wsIndi 130 1 "2028-01-01"
$wsIndi.Cells.Item(130, 2).Value = 1301.54878945716
$wsIndi.Cells.Item(130, 3).Value = 178.375530074545
Set-DateText $wsIndi 131 1 "2028-04-01"
$wsIndi.Cells.Item(131, 2).Value = 1307.94113123125
$wsIndi.Cells.Item(131, 3).Value = 179.251592010608
Set-DateText $wsIndi 132 1 "2028-07-01"
$wsIndi.Cells.Item(132, 2).Value = 1318.97608024656
$wsIndi.Cells.Item(132, 3).Value = 180.763917092767
Set-DateText $wsIndi 133 1 "2028-10-01"
$wsIndi.Cells.Item(133, 2).Value = 1329.66733753017
$wsIndi.Cells.Item(133, 3).Value = 182.22913968033
Set-DateText $wsIndi 134 1 "2029-01-01"
$wsIndi.Cells.Item(134, 2).Value = 1333.92697813096
$wsIndi.Cells.Item(134, 3).Value = 182.812917757839
Set-DateText $wsIndi 135 1 "2029-04-01"
$wsIndi.Cells.Item(135, 2).Value = 1340.84484405888
$wsIndi.Cells.Item(135, 3).Value = 183.76100208005
Set-DateText $wsIndi 136 1 "2029-07-01"
$wsIndi.Cells.Item(136, 2).Value = 1350.77704720087
$wsIndi.Cells.Item(136, 3).Value = 185.122197307314
Set-DateText $wsIndi 137 1 "2029-10-01"
$wsIndi.Cells.Item(137, 2).Value = 1360.53536447346
$wsIndi.Cells.Item(137, 3).Value = 186.459561707507

# ---- ETALONNAGE: refresh recalculated rows 27-31 ----
$wsEtal = $wb.Worksheets.Item("ETALONNAGE")
$wsEtal.Cells.Item(27, 2).Value = 2696241.20362179
$wsEtal.Cells.Item(27, 4).Value = 9.43239781472145
$wsEtal.Cells.Item(28, 2).Value = 2825901.05035834
$wsEtal.Cells.Item(28, 4).Value = 4.80891125624714
$wsEtal.Cells.Item(29, 2).Value = 2940787.7456802
$wsEtal.Cells.Item(29, 4).Value = 4.06548896350398
$wsEtal.Cells.Item(30, 2).Value = 3035886.54632515
$wsEtal.Cells.Item(30, 4).Value = 3.23378662008649
$wsEtal.Cells.Item(31, 2).Value = 3126352.91827335
$wsEtal.Cells.Item(31, 4).Value = 2.97989962957299

# ---- ETALONNAGE: append forecast rows 32-35 (annual 2026-2029) ----
Set-DateText $wsEtal 32 1 "2026-01-01"
$wsEtal.Cells.Item(32, 2).Value = 3212022.98159255
$wsEtal.Cells.Item(32, 3).Value = 171.235897705177
$wsEtal.Cells.Item(32, 4).Value = 2.74025567678144
$wsEtal.Cells.Item(32, 5).Value = 2.74090844496171
$wsEtal.Cells.Item(32, 6).Value = "Acceptable"
Set-DateText $wsEtal 33 1 "2027-01-01"
$wsEtal.Cells.Item(33, 2).Value = 3295847.27778426
$wsEtal.Cells.Item(33, 3).Value = 175.731467625905
$wsEtal.Cells.Item(33, 4).Value = 2.60970412329191
$wsEtal.Cells.Item(33, 5).Value = 2.62536651541818
$wsEtal.Cells.Item(33, 6).Value = "Acceptable"
Set-DateText $wsEtal 34 1 "2028-01-01"
$wsEtal.Cells.Item(34, 2).Value = 3378140.17749033
$wsEtal.Cells.Item(34, 3).Value = 180.155044714563
$wsEtal.Cells.Item(34, 4).Value = 2.4968662917352
$wsEtal.Cells.Item(34, 5).Value = 2.51723675242663
$wsEtal.Cells.Item(34, 6).Value = "Acceptable"
Set-DateText $wsEtal 35 1 "2029-01-01"
$wsEtal.Cells.Item(35, 2).Value = 3459619.70742599
$wsEtal.Cells.Item(35, 3).Value = 184.538919713177
$wsEtal.Cells.Item(35, 4).Value = 2.41196414756806
$wsEtal.Cells.Item(35, 5).Value = 2.43339008661148
$wsEtal.Cells.Item(35, 6).Value = "Acceptable"

# ---- PREVISION: refresh recalculated rows 98-117 ----
$wsPrev = $wb.Worksheets.Item("PREVISION")
$wsPrev.Cells.Item(98, 2).Value = 649904.449852078
$wsPrev.Cells.Item(98, 4).Value = -41526.917546703
$wsPrev.Cells.Item(99, 2).Value = 663964.757653398
$wsPrev.Cells.Item(99, 4).Value = -39566.2709986051
$wsPrev.Cells.Item(100, 2).Value = 676117.192732277
$wsPrev.Cells.Item(100, 4).Value = -35644.9779024093
$wsPrev.Cells.Item(101, 2).Value = 706254.803384042
$wsPrev.Cells.Item(101, 4).Value = -29763.0382581157
$wsPrev.Cells.Item(102, 2).Value = 689024.468353722
$wsPrev.Cells.Item(102, 4).Value = -21920.4520657242
$wsPrev.Cells.Item(103, 2).Value = 678763.594581213
$wsPrev.Cells.Item(103, 4).Value = -15743.2762075106
$wsPrev.Cells.Item(104, 2).Value = 715825.813389977
$wsPrev.Cells.Item(104, 4).Value = -11231.5106834749
$wsPrev.Cells.Item(105, 2).Value = 742287.174033425
$wsPrev.Cells.Item(105, 4).Value = -8385.15549361713
$wsPrev.Cells.Item(106, 2).Value = 721123.351092959
$wsPrev.Cells.Item(106, 4).Value = -7204.21063793729
$wsPrev.Cells.Item(107, 2).Value = 723961.944331744
$wsPrev.Cells.Item(107, 4).Value = -6089.75022053941
$wsPrev.Cells.Item(108, 2).Value = 739236.470084766
$wsPrev.Cells.Item(108, 4).Value = -5041.77424142352
$wsPrev.Cells.Item(109, 2).Value = 756465.980170728
$wsPrev.Cells.Item(109, 4).Value = -4060.2827005896
$wsPrev.Cells.Item(110, 2).Value = 747562.126712389
$wsPrev.Cells.Item(110, 4).Value = -3145.27559803766
$wsPrev.Cells.Item(111, 2).Value = 747289.723590938
$wsPrev.Cells.Item(111, 4).Value = -2396.82271206747
$wsPrev.Cells.Item(112, 2).Value = 763523.565420536
$wsPrev.Cells.Item(112, 4).Value = -1814.92404267904
$wsPrev.Cells.Item(113, 2).Value = 777511.130601284
$wsPrev.Cells.Item(113, 4).Value = -1399.57958987236
$wsPrev.Cells.Item(114, 2).Value = 771383.253340462
$wsPrev.Cells.Item(114, 4).Value = -1150.78935364746
$wsPrev.Cells.Item(115, 2).Value = 774218.252053636
$wsPrev.Cells.Item(115, 4).Value = -933.331230137956
$wsPrev.Cells.Item(116, 2).Value = 785013.022875863
$wsPrev.Cells.Item(116, 4).Value = -747.205219343869
$wsPrev.Cells.Item(117, 2).Value = 795738.390003386
$wsPrev.Cells.Item(117, 4).Value = -592.411321265194

# ---- PREVISION: append forecast rows 118-133 (2026Q1 .. 2029Q4) ----
Set-DateText $wsPrev 118 1 "2026-01-01"
$wsPrev.Cells.Item(118, 2).Value = 793784.095683984
$wsPrev.Cells.Item(118, 3).Value = 42.3153481549725
$wsPrev.Cells.Item(118, 4).Value = -468.949535901935
Set-DateText $wsPrev 119 1 "2026-04-01"
$wsPrev.Cells.Item(119, 2).Value = 796672.995455177
$wsPrev.Cells.Item(119, 3).Value = 42.4687036026278
$wsPrev.Cells.Item(119, 4).Value = -365.8494444654
Set-DateText $wsPrev 120 1 "2026-07-01"
$wsPrev.Cells.Item(120, 2).Value = 806286.718209278
$wsPrev.Cells.Item(120, 3).Value = 42.9845224401467
$wsPrev.Cells.Item(120, 4).Value = -283.111046955584
Set-DateText $wsPrev 121 1 "2026-10-01"
$wsPrev.Cells.Item(121, 2).Value = 815279.172244115
$wsPrev.Cells.Item(121, 3).Value = 43.4673235074302
$wsPrev.Cells.Item(121, 4).Value = -220.734343372483
Set-DateText $wsPrev 122 1 "2027-01-01"
$wsPrev.Cells.Item(122, 2).Value = 815286.320888311
$wsPrev.Cells.Item(122, 3).Value = 43.4668239792683
$wsPrev.Cells.Item(122, 4).Value = -178.719333716109
Set-DateText $wsPrev 123 1 "2027-04-01"
$wsPrev.Cells.Item(123, 2).Value = 819084.330351387
$wsPrev.Cells.Item(123, 3).Value = 43.6705352442817
$wsPrev.Cells.Item(123, 4).Value = -142.765150663712
Set-DateText $wsPrev 124 1 "2027-07-01"
$wsPrev.Cells.Item(124, 2).Value = 826915.889095407
$wsPrev.Cells.Item(124, 3).Value = 44.0915225783221
$wsPrev.Cells.Item(124, 4).Value = -112.871794215288
Set-DateText $wsPrev 125 1 "2027-10-01"
$wsPrev.Cells.Item(125, 2).Value = 834560.737449155
$wsPrev.Cells.Item(125, 3).Value = 44.5025858240325
$wsPrev.Cells.Item(125, 4).Value = -89.0392643708346
Set-DateText $wsPrev 126 1 "2028-01-01"
$wsPrev.Cells.Item(126, 2).Value = 836263.528826173
$wsPrev.Cells.Item(126, 3).Value = 44.5938825186363
$wsPrev.Cells.Item(126, 4).Value = -71.2675611303593
Set-DateText $wsPrev 127 1 "2028-04-01"
$wsPrev.Cells.Item(127, 2).Value = 840337.68166146
$wsPrev.Cells.Item(127, 3).Value = 44.8128980026519
$wsPrev.Cells.Item(127, 4).Value = -56.1160578538458
Set-DateText $wsPrev 128 1 "2028-07-01"
$wsPrev.Cells.Item(128, 2).Value = 847365.46873716
$wsPrev.Cells.Item(128, 3).Value = 45.1909792731918
$wsPrev.Cells.Item(128, 4).Value = -43.5847545412982
Set-DateText $wsPrev 129 1 "2028-10-01"
$wsPrev.Cells.Item(129, 2).Value = 854173.498265534
$wsPrev.Cells.Item(129, 3).Value = 45.5572849200826
$wsPrev.Cells.Item(129, 4).Value = -33.6736511927158
Set-DateText $wsPrev 130 1 "2029-01-01"
$wsPrev.Cells.Item(130, 2).Value = 856887.279136187
$wsPrev.Cells.Item(130, 3).Value = 45.7032294394596
$wsPrev.Cells.Item(130, 4).Value = -26.3827478080968
Set-DateText $wsPrev 131 1 "2029-04-01"
$wsPrev.Cells.Item(131, 2).Value = 861292.101178299
$wsPrev.Cells.Item(131, 3).Value = 45.9402505200124
$wsPrev.Cells.Item(131, 4).Value = -20.9145702696316
Set-DateText $wsPrev 132 1 "2029-07-01"
$wsPrev.Cells.Item(132, 2).Value = 867614.601932794
$wsPrev.Cells.Item(132, 3).Value = 46.2805493268285
$wsPrev.Cells.Item(132, 4).Value = -17.2691185773221
Set-DateText $wsPrev 133 1 "2029-10-01"
$wsPrev.Cells.Item(133, 2).Value = 873825.725178706
$wsPrev.Cells.Item(133, 3).Value = 46.6148904268767
$wsPrev.Cells.Item(133, 4).Value = -15.4463927311674

# ---- VATRIM: refresh recalculated rows 98-117 ----
$wsVat = $wb.Worksheets.Item("VATRIM")
$wsVat.Cells.Item(98, 2).Value = 649904.449852078
$wsVat.Cells.Item(99, 2).Value = 663964.757653398
$wsVat.Cells.Item(100, 2).Value = 676117.192732277
$wsVat.Cells.Item(101, 2).Value = 706254.803384042
$wsVat.Cells.Item(102, 2).Value = 689024.468353722
$wsVat.Cells.Item(103, 2).Value = 678763.594581213
$wsVat.Cells.Item(104, 2).Value = 715825.813389977
$wsVat.Cells.Item(105, 2).Value = 742287.174033425
$wsVat.Cells.Item(106, 2).Value = 721123.351092959
$wsVat.Cells.Item(107, 2).Value = 723961.944331744
$wsVat.Cells.Item(108, 2).Value = 739236.470084766
$wsVat.Cells.Item(109, 2).Value = 756465.980170728
$wsVat.Cells.Item(110, 2).Value = 747562.126712389
$wsVat.Cells.Item(111, 2).Value = 747289.723590938
$wsVat.Cells.Item(112, 2).Value = 763523.565420536
$wsVat.Cells.Item(113, 2).Value = 777511.130601284
$wsVat.Cells.Item(114, 2).Value = 771383.253340462
$wsVat.Cells.Item(115, 2).Value = 774218.252053636
$wsVat.Cells.Item(116, 2).Value = 785013.022875863
$wsVat.Cells.Item(117, 2).Value = 795738.390003386

# ---- VATRIM: append forecast rows 118-133 (2026Q1 .. 2029Q4) ----
Set-DateText $wsVat 118 1 "2026-01-01"
$wsVat.Cells.Item(118, 2).Value = 793784.095683984
Set-DateText $wsVat 119 1 "2026-04-01"
$wsVat.Cells.Item(119, 2).Value = 796672.995455177
Set-DateText $wsVat 120 1 "2026-07-01"
$wsVat.Cells.Item(120, 2).Value = 806286.718209278
Set-DateText $wsVat 121 1 "2026-10-01"
$wsVat.Cells.Item(121, 2).Value = 815279.172244115
Set-DateText $wsVat 122 1 "2027-01-01"
$wsVat.Cells.Item(122, 2).Value = 815286.320888311
Set-DateText $wsVat 123 1 "2027-04-01"
$wsVat.Cells.Item(123, 2).Value = 819084.330351387
Set-DateText $wsVat 124 1 "2027-07-01"
$wsVat.Cells.Item(124, 2).Value = 826915.889095407
Set-DateText $wsVat 125 1 "2027-10-01"
$wsVat.Cells.Item(125, 2).Value = 834560.737449155
Set-DateText $wsVat 126 1 "2028-01-01"
$wsVat.Cells.Item(126, 2).Value = 836263.528826173
Set-DateText $wsVat 127 1 "2028-04-01"
$wsVat.Cells.Item(127, 2).Value = 840337.68166146
Set-DateText $wsVat 128 1 "2028-07-01"
$wsVat.Cells.Item(128, 2).Value = 847365.46873716
Set-DateText $wsVat 129 1 "2028-10-01"
$wsVat.Cells.Item(129, 2).Value = 854173.498265534
Set-DateText $wsVat 130 1 "2029-01-01"
$wsVat.Cells.Item(130, 2).Value = 856887.279136187
Set-DateText $wsVat 131 1 "2029-04-01"
$wsVat.Cells.Item(131, 2).Value = 861292.101178299
Set-DateText $wsVat 132 1 "2029-07-01"
$wsVat.Cells.Item(132, 2).Value = 867614.601932794
Set-DateText $wsVat 133 1 "2029-10-01"
$wsVat.Cells.Item(133, 2).Value = 873825.725178706

